# Converts an 8-bit R,G,B triple into the OLE/COM "long" RGB representation
# PowerPoint's object model expects (0x00BBGGRR).
function RGBVal($r, $g, $b) {
    return $b * 65536 + $g * 256 + $r
}

$p = $ppt.ActivePresentation

# --- 1) Table style ids: 3 tables (slides 14, 15 and 16) move from the
#     deck's custom "Table_0" style to the built-in style
#     {81B09CB9-A143-4304-B9E3-818962A0AE93}.
$newTableStyle = "{81B09CB9-A143-4304-B9E3-818962A0AE93}"

$s14 = $p.Slides.Item(14)
$s14.Shapes.Item(1).Table.ApplyStyle($newTableStyle)

$s15 = $p.Slides.Item(15)
$s15.Shapes.Item(1).Table.ApplyStyle($newTableStyle)

$s16 = $p.Slides.Item(16)
$s16.Shapes.Item(1).Table.ApplyStyle($newTableStyle)

# --- 2) Theme: the deck's design swaps from the custom "Integral" / "Red
#     Violet" colour scheme back to the default "Office Theme" / "Office"
#     colour scheme. The font scheme and format scheme are identical
#     between the two themes, so only the 12 theme colours need updating.
$tcs = $s14.ThemeColorScheme

$tcs.Item(1).RGB  = RGBVal 0x00 0x00 0x00   # dk1
$tcs.Item(2).RGB  = RGBVal 0xFF 0xFF 0xFF   # lt1
$tcs.Item(3).RGB  = RGBVal 0x44 0x54 0x6A   # dk2
$tcs.Item(4).RGB  = RGBVal 0xE7 0xE6 0xE6   # lt2
$tcs.Item(5).RGB  = RGBVal 0x5B 0x9B 0xD5   # accent1
$tcs.Item(6).RGB  = RGBVal 0xED 0x7D 0x31   # accent2
$tcs.Item(7).RGB  = RGBVal 0xA5 0xA5 0xA5   # accent3
$tcs.Item(8).RGB  = RGBVal 0xFF 0xC0 0x00   # accent4
$tcs.Item(9).RGB  = RGBVal 0x44 0x72 0xC4   # accent5
$tcs.Item(10).RGB = RGBVal 0x70 0xAD 0x47   # accent6
$tcs.Item(11).RGB = RGBVal 0x05 0x63 0xC1   # hlink
$tcs.Item(12).RGB = RGBVal 0x95 0x4F 0x72   # folHlink
